# Dashboard and all updates
# Appends the new income/expense rows (3-7) to the Income sheet, mirroring
# the date formatting already used in row 2 (column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateValue = 45817.54185799768

$rows = @(
    @("Salaryy",       1500),
    @("Sdhjbsfryy",    15000),
    @("groceries",     15000),
    @("maid",          15000),
    @("jygvhcfxcgjv",  1500)
)

$r = 3
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $dateValue
    # Same built-in date format (numFmtId 14) already applied to C2.
    $ws.Cells.Item($r, 3).NumberFormat = "mm-dd-yy"
    $r++
}
